$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = 'maa://24702 (94.31), maa://25390 (95.98), maa://36681 (87.01)'
$ws.Range('L2').Value = '*maa://24633 (56.05), *maa://30515 (69.31), *maa://34787 (73.61), ***maa://20792 (11.93), maa://39402 (89.36), ***maa://29083 (27.78)'
$ws.Range('T2').Value = 'maa://22742 (91.08), *maa://20791 (63.89)'
$ws.Range('AB2').Value = 'maa://21246 (91.38), maa://36684 (94.79), ***maa://22731 (6.67)'
$ws.Range('AF2').Value = 'maa://25251 (92.08), ***maa://21730 (22.86), ***maa://39501 (18.18), *maa://36675 (60.0)'
$ws.Range('H3').Value = 'maa://21247 (98.47), *maa://22748 (60.0)'
$ws.Range('L3').Value = '*maa://22880 (66.13), maa://20276 (85.19), *maa://22749 (72.73)'
$ws.Range('P3').Value = 'maa://21249 (94.22), maa://26254 (96.3)'
$ws.Range('T3').Value = 'maa://24617 (89.09), **maa://20790 (44.12), ***maa://37170 (18.97)'
$ws.Range('T4').Value = 'maa://32509 (98.04), maa://27295 (84.13), maa://22754 (90.41), *maa://21746 (55.81), *maa://31008 (78.57)'
$ws.Range('X4').Value = '**maa://32495 (47.91), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (90.91)'
$ws.Range('A8').Value = '更新日期：2025.01.09 13:18:22'
$ws.Range('X8').Value = 'maa://21411 (95.97)'
$ws.Range('AF10').Value = '*maa://25021 (55.29), *maa://22733 (59.38), maa://22761 (100.0)'
$ws.Range('X11').Value = 'maa://36713 (97.96)'
$ws.Range('D13').Value = 'maa://24999 (91.83), maa://36673 (92.75), maa://25001 (85.51)'
$ws.Range('X13').Value = 'maa://34957 (80.88), *maa://22768 (51.61)'
$ws.Range('AF13').Value = '**maa://22737 (32.37), maa://39883 (91.23), *maa://39885 (56.0)'
$ws.Range('L14').Value = 'maa://26245 (96.55), maa://21288 (96.3), maa://39841 (95.35), maa://36682 (97.37)'
$ws.Range('P14').Value = 'maa://23250 (98.63), maa://20107 (87.1), maa://22772 (100.0), **maa://22745 (50.0)'
$ws.Range('D15').Value = '*maa://22743 (77.55), maa://22734 (84.03), *maa://30808 (65.08), **maa://36048 (34.04), maa://45058 (100.0)'
$ws.Range('H15').Value = 'maa://24304 (88.24), maa://21478 (91.43)'
$ws.Range('AF15').Value = 'maa://21364 (81.09), *maa://22766 (70.27), *maa://36666 (78.16)'
$ws.Range('T16').Value = 'maa://22729 (94.77), *maa://28648 (68.25), maa://36674 (83.33)'
$ws.Range('X16').Value = 'maa://28501 (97.85), maa://28051 (96.0)'
$ws.Range('H17').Value = 'maa://22430 (88.77), maa://39599 (86.84)'
$ws.Range('D18').Value = 'maa://24570 (97.16)'
$ws.Range('H18').Value = 'maa://24421 (90.16)'
$ws.Range('T19').Value = 'maa://24386 (99.06)'
$ws.Range('D20').Value = 'maa://21432 (90.07), maa://25198 (93.14), *maa://20795 (51.18), maa://36680 (93.55)'
$ws.Range('AB21').Value = 'maa://21443 (80.06), ***maa://23820 (29.31)'
$ws.Range('H22').Value = 'maa://25236 (96.39), **maa://21678 (48.94), **maa://22735 (42.86)'
$ws.Range('L22').Value = 'maa://27127 (84.91), *maa://22751 (73.85)'
$ws.Range('L23').Value = 'maa://39756 (94.53), maa://39875 (93.85)'
$ws.Range('X24').Value = 'maa://29988 (86.46), maa://23504 (93.18), **maa://22892 (39.31), *maa://25141 (76.56), *maa://36663 (78.57), ***maa://22815 (23.08)'
$ws.Range('AF25').Value = 'maa://20108 (96.27), maa://24621 (96.58), maa://36676 (96.88), maa://22771 (85.71), **maa://37772 (50.0)'
$ws.Range('X28').Value = 'maa://39929 (89.97), maa://41749 (91.94), ***maa://39723 (14.29)'
$ws.Range('AB30').Value = 'maa://42979 (96.67), maa://45045 (100.0)'
$ws.Range('L31').Value = 'maa://35926 (93.38), maa://36258 (84.0), *maa://43904 (72.73)'
$ws.Range('H32').Value = 'maa://21895 (97.38), maa://36667 (98.41), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range('T32').Value = 'maa://42859 (96.63), maa://41108 (87.76), maa://41238 (96.43)'
$ws.Range('L35').Value = 'maa://41296 (96.83)'
$ws.Range('H39').Value = 'maa://25199 (84.82), maa://36670 (87.91), maa://30434 (89.55), ***maa://25036 (16.0), *maa://44165 (66.67), *maa://45059 (75.0)'
$ws.Range('H45').Value = 'maa://21229 (84.86), maa://30807 (95.52), *maa://22767 (55.0), ***maa://20796 (13.79), *maa://42459 (80.0)'
$ws.Range('T45').Value = '**maa://39364 (39.13)'
$ws.Range('H46').Value = 'maa://35931 (92.41), maa://43901 (90.0)'
$ws.Range('H47').Value = 'maa://27410 (96.24), maa://29661 (97.86), maa://28038 (84.62)'
$ws.Range('H55').Value = 'maa://32532 (92.02)'
$ws.Range('H60').Value = '*maa://40438 (61.7)'
$ws.Range('H62').Value = 'maa://42981 (96.97), maa://43903 (100.0)'
